$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
 MATCH (f:file)-->(parent)
MATCH (c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c)
 WHERE  s.clinical_study_designation in ['NCATS-COP01'] and s.clinical_study_type in ['Transcriptomics']  and samp.sample_site in ['Mandible, Mucosa'] and labels(parent)[0] IN ['study']
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
 MATCH (f:file)-->(parent)
WHERE s.clinical_study_designation in ['NCATS-COP01'] and s.clinical_study_type in ['Transcriptomics'] and samp.sample_site in ['Mandible, Mucosa'] and labels(parent)[0] IN ['study']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (r:registration)-->(c)
MATCH (f)-->(samp:sample)
WHERE s.clinical_study_designation in ['NCATS-COP01'] and s.clinical_study_type in ['Transcriptomics'] and samp.sample_site in ['Mandible, Mucosa'] and labels(parent)[0] IN ['study'] 
 MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@

$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (samp:sample)-->(c)
MATCH (c)<--(demo:demographic)
MATCH (f:file)-->(parent)
WHERE s.clinical_study_designation in ['NCATS-COP01'] and s.clinical_study_type in ['Transcriptomics'] and samp.sample_site in ['Mandible, Mucosa'] and labels(parent)[0] IN ['study'] 
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$statQuery = @'
MATCH (s:study)<--(c:case)
WHERE s.clinical_study_designation IN ['NCATS-COP01']
OPTIONAL MATCH (cf:file)-->(c)
OPTIONAL MATCH (sf:file)-->(s)
 MATCH (samp:sample)-->(c)
 MATCH (f:file)-->(parent)
where s.clinical_study_type in ['Transcriptomics'] and samp.sample_site in ['Mandible, Mucosa'] and labels(parent)[0] IN ['study']
OPTIONAL MATCH (s)-->(p:program)
OPTIONAL MATCH (cf:file)-[*]->(c)
OPTIONAL MATCH (p:program)-->(s)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_WebData.xlsx"

$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_WebData.xlsx"

$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_WebData.xlsx"

$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC21_Canine_StudyNCATSCOP01-StudyType_SampleSite_FileAsso_WebData.xlsx"

$ws.Rows.Item(2).RowHeight = 251.25
$ws.Rows.Item(3).RowHeight = 177
$ws.Rows.Item(4).RowHeight = 164.25
$ws.Rows.Item(5).RowHeight = 204

$ws.Range("C5").Select() | Out-Null
